$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "No. of R3 Excel Row's to Execute" value from 200 to 100
$ws.Range("D2").Value = "100"

# Move the active selection to D2
$ws.Range("D2").Select()
